# Auto-generated edit script: refresh market-price-driven columns
# (currentAveragePrice / NQ / HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 3000  # H40: 0 -> 3000
$ws.Cells.Item(40, 9).Value = 3000  # I40: 0 -> 3000
$ws.Cells.Item(40, 11).Value = 3000  # K40: 0 -> 3000
$ws.Cells.Item(40, 13).Value = -2825  # M40: None -> -2825
$ws.Cells.Item(58, 8).Value = 35714696  # H58: 25000302 -> 35714696
$ws.Cells.Item(58, 9).Value = 35714696  # I58: 25000302 -> 35714696
$ws.Cells.Item(58, 11).Value = 107144088  # K58: 75000906 -> 107144088
$ws.Cells.Item(58, 13).Value = -107143938  # M58: -75000756 -> -107143938
$ws.Cells.Item(74, 8).Value = 2008.8182  # H74: 2749.8333 -> 2008.8182
$ws.Cells.Item(74, 9).Value = 2008.8182  # I74: 2749.8333 -> 2008.8182
$ws.Cells.Item(74, 11).Value = 2008.8182  # K74: 2749.8333 -> 2008.8182
$ws.Cells.Item(74, 13).Value = -1072.8182  # M74: -1813.8333 -> -1072.8182
$ws.Cells.Item(77, 8).Value = 2008.8182  # H77: 2749.8333 -> 2008.8182
$ws.Cells.Item(77, 9).Value = 2008.8182  # I77: 2749.8333 -> 2008.8182
$ws.Cells.Item(77, 11).Value = 10044.091  # K77: 13749.1665 -> 10044.091
$ws.Cells.Item(77, 13).Value = -5364.091  # M77: -9069.166499999999 -> -5364.091
$ws.Cells.Item(129, 8).Value = 2808.8572  # H129: 2788.48 -> 2808.8572
$ws.Cells.Item(129, 10).Value = 2987.0417  # J129: 2988.238 -> 2987.0417
$ws.Cells.Item(129, 12).Value = 8961.125100000001  # L129: 8964.714 -> 8961.125100000001
$ws.Cells.Item(129, 14).Value = -18961.1251  # N129: -18964.714 -> -18961.1251
$ws.Cells.Item(131, 8).Value = 4462.067  # H131: 4735.9375 -> 4462.067
$ws.Cells.Item(131, 9).Value = 448.55554  # I131: 431.33334 -> 448.55554
$ws.Cells.Item(131, 10).Value = 10482.333  # J131: 10270.429 -> 10482.333
$ws.Cells.Item(131, 11).Value = 1345.66662  # K131: 1294.00002 -> 1345.66662
$ws.Cells.Item(131, 12).Value = 31446.999  # L131: 30811.287 -> 31446.999
$ws.Cells.Item(131, 13).Value = 3694.33338  # M131: 3745.99998 -> 3694.33338
$ws.Cells.Item(131, 14).Value = -41526.999  # N131: -40891.287 -> -41526.999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3331.1428  # H32: 3295.6 -> 3331.1428
$ws.Cells.Item(32, 9).Value = 2660.2444  # I32: 2661.4666 -> 2660.2444
$ws.Cells.Item(32, 10).Value = 10878.75  # J32: 9002.799999999999 -> 10878.75
$ws.Cells.Item(32, 11).Value = 2660.2444  # K32: 2661.4666 -> 2660.2444
$ws.Cells.Item(32, 12).Value = 10878.75  # L32: 9002.799999999999 -> 10878.75
$ws.Cells.Item(32, 13).Value = -2373.2444  # M32: -2374.4666 -> -2373.2444
$ws.Cells.Item(32, 14).Value = -11452.75  # N32: -9576.799999999999 -> -11452.75
$ws.Cells.Item(45, 8).Value = 1385.1538  # H45: 1539.8572 -> 1385.1538
$ws.Cells.Item(45, 9).Value = 891.5454999999999  # I45: 1050.8182 -> 891.5454999999999
$ws.Cells.Item(45, 10).Value = 4100  # J45: 3333 -> 4100
$ws.Cells.Item(45, 11).Value = 891.5454999999999  # K45: 1050.8182 -> 891.5454999999999
$ws.Cells.Item(45, 12).Value = 4100  # L45: 3333 -> 4100
$ws.Cells.Item(45, 13).Value = -514.5454999999999  # M45: -673.8181999999999 -> -514.5454999999999
$ws.Cells.Item(45, 14).Value = -4854  # N45: -4087 -> -4854
$ws.Cells.Item(88, 8).Value = 1250.8823  # H88: 1228.4 -> 1250.8823
$ws.Cells.Item(88, 9).Value = 1262.1428  # I88: 1037.8889 -> 1262.1428
$ws.Cells.Item(88, 10).Value = 1243  # J88: 1335.5625 -> 1243
$ws.Cells.Item(88, 11).Value = 1262.1428  # K88: 1037.8889 -> 1262.1428
$ws.Cells.Item(88, 12).Value = 1243  # L88: 1335.5625 -> 1243
$ws.Cells.Item(88, 13).Value = -856.1428000000001  # M88: -631.8888999999999 -> -856.1428000000001
$ws.Cells.Item(88, 14).Value = -2055  # N88: -2147.5625 -> -2055
$ws.Cells.Item(91, 8).Value = 1250.8823  # H91: 1228.4 -> 1250.8823
$ws.Cells.Item(91, 9).Value = 1262.1428  # I91: 1037.8889 -> 1262.1428
$ws.Cells.Item(91, 10).Value = 1243  # J91: 1335.5625 -> 1243
$ws.Cells.Item(91, 11).Value = 1262.1428  # K91: 1037.8889 -> 1262.1428
$ws.Cells.Item(91, 12).Value = 1243  # L91: 1335.5625 -> 1243
$ws.Cells.Item(91, 13).Value = 141.8571999999999  # M91: 366.1111000000001 -> 141.8571999999999
$ws.Cells.Item(91, 14).Value = -4051  # N91: -4143.5625 -> -4051
$ws.Cells.Item(92, 8).Value = 35234.54  # H92: 39304.9 -> 35234.54
$ws.Cells.Item(92, 10).Value = 35234.54  # J92: 39304.9 -> 35234.54
$ws.Cells.Item(92, 12).Value = 35234.54  # L92: 39304.9 -> 35234.54
$ws.Cells.Item(92, 14).Value = -40226.54  # N92: -44296.9 -> -40226.54
$ws.Cells.Item(108, 8).Value = 30000  # H108: 26810.5 -> 30000
$ws.Cells.Item(108, 9).Value = 0  # I108: 26810.5 -> 0
$ws.Cells.Item(108, 10).Value = 30000  # J108: 0 -> 30000
$ws.Cells.Item(108, 11).Value = 0  # K108: 26810.5 -> 0
$ws.Cells.Item(108, 12).Value = 30000  # L108: 0 -> 30000
$ws.Cells.Item(108, 13).ClearContents()  # M108: remove (was -22970.5)
$ws.Cells.Item(108, 14).Value = -37680  # N108: None -> -37680

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 30531618  # H86: 28835362 -> 30531618
$ws.Cells.Item(86, 9).Value = 1813.875  # I86: 1787.4286 -> 1813.875
$ws.Cells.Item(86, 10).Value = 57669224  # J86: 47184000 -> 57669224
$ws.Cells.Item(86, 11).Value = 1813.875  # K86: 1787.4286 -> 1813.875
$ws.Cells.Item(86, 12).Value = 57669224  # L86: 47184000 -> 57669224
$ws.Cells.Item(86, 13).Value = -690.875  # M86: -664.4286 -> -690.875
$ws.Cells.Item(86, 14).Value = -57671470  # N86: -47186246 -> -57671470
$ws.Cells.Item(89, 8).Value = 30531618  # H89: 28835362 -> 30531618
$ws.Cells.Item(89, 9).Value = 1813.875  # I89: 1787.4286 -> 1813.875
$ws.Cells.Item(89, 10).Value = 57669224  # J89: 47184000 -> 57669224
$ws.Cells.Item(89, 11).Value = 9069.375  # K89: 8937.143 -> 9069.375
$ws.Cells.Item(89, 12).Value = 288346120  # L89: 235920000 -> 288346120
$ws.Cells.Item(89, 13).Value = -3453.375  # M89: -3321.143 -> -3453.375
$ws.Cells.Item(89, 14).Value = -288357352  # N89: -235931232 -> -288357352
$ws.Cells.Item(99, 8).Value = 3432.75  # H99: 3177.1538 -> 3432.75
$ws.Cells.Item(99, 9).Value = 1854.7778  # I99: 1680.3 -> 1854.7778
$ws.Cells.Item(99, 11).Value = 1854.7778  # K99: 1680.3 -> 1854.7778
$ws.Cells.Item(99, 13).Value = -356.7778000000001  # M99: -182.3 -> -356.7778000000001
$ws.Cells.Item(107, 8).Value = 5479.9653  # H107: 5154.9062 -> 5479.9653
$ws.Cells.Item(107, 9).Value = 4051  # I107: 3785.1304 -> 4051
$ws.Cells.Item(107, 11).Value = 4051  # K107: 3785.1304 -> 4051
$ws.Cells.Item(107, 13).Value = -2131  # M107: -1865.1304 -> -2131
$ws.Cells.Item(125, 8).Value = 75000  # H125: 0 -> 75000
$ws.Cells.Item(125, 10).Value = 75000  # J125: 0 -> 75000
$ws.Cells.Item(125, 12).Value = 75000  # L125: 0 -> 75000
$ws.Cells.Item(125, 14).Value = -84840  # N125: None -> -84840

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1521.1765  # H16: 1533.5883 -> 1521.1765
$ws.Cells.Item(16, 9).Value = 1263.3334  # I16: 1280.9166 -> 1263.3334
$ws.Cells.Item(16, 11).Value = 1263.3334  # K16: 1280.9166 -> 1263.3334
$ws.Cells.Item(16, 13).Value = -976.3334  # M16: -993.9166 -> -976.3334
$ws.Cells.Item(31, 8).Value = 6620.9585  # H31: 5680.5864 -> 6620.9585
$ws.Cells.Item(31, 9).Value = 2414.1428  # I31: 1894.4166 -> 2414.1428
$ws.Cells.Item(31, 11).Value = 2414.1428  # K31: 1894.4166 -> 2414.1428
$ws.Cells.Item(31, 13).Value = -2119.1428  # M31: -1599.4166 -> -2119.1428
$ws.Cells.Item(34, 8).Value = 6620.9585  # H34: 5680.5864 -> 6620.9585
$ws.Cells.Item(34, 9).Value = 2414.1428  # I34: 1894.4166 -> 2414.1428
$ws.Cells.Item(34, 11).Value = 2414.1428  # K34: 1894.4166 -> 2414.1428
$ws.Cells.Item(34, 13).Value = -2212.1428  # M34: -1692.4166 -> -2212.1428
$ws.Cells.Item(113, 8).Value = 1521.1765  # H113: 1533.5883 -> 1521.1765
$ws.Cells.Item(113, 9).Value = 1263.3334  # I113: 1280.9166 -> 1263.3334
$ws.Cells.Item(113, 11).Value = 1263.3334  # K113: 1280.9166 -> 1263.3334
$ws.Cells.Item(113, 13).Value = 906.6666  # M113: 889.0834 -> 906.6666
$ws.Cells.Item(137, 8).Value = 127245.37  # H137: 137523.22 -> 127245.37
$ws.Cells.Item(137, 10).Value = 89974.89999999999  # J137: 92219.875 -> 89974.89999999999
$ws.Cells.Item(137, 12).Value = 89974.89999999999  # L137: 92219.875 -> 89974.89999999999
$ws.Cells.Item(137, 14).Value = -100174.9  # N137: -102419.875 -> -100174.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 64.8  # H2: 57.565216 -> 64.8
$ws.Cells.Item(2, 9).Value = 39.25  # I2: 33.266666 -> 39.25
$ws.Cells.Item(2, 11).Value = 235.5  # K2: 199.599996 -> 235.5
$ws.Cells.Item(2, 13).Value = -122.5  # M2: -86.599996 -> -122.5
$ws.Cells.Item(5, 8).Value = 781.8889  # H5: 793.3 -> 781.8889
$ws.Cells.Item(5, 10).Value = 999.6667  # J5: 973.75 -> 999.6667
$ws.Cells.Item(5, 12).Value = 2999.0001  # L5: 2921.25 -> 2999.0001
$ws.Cells.Item(5, 14).Value = -3223.0001  # N5: -3145.25 -> -3223.0001
$ws.Cells.Item(33, 8).Value = 207.92308  # H33: 300.85 -> 207.92308
$ws.Cells.Item(33, 9).Value = 217.875  # I33: 168.81818 -> 217.875
$ws.Cells.Item(33, 10).Value = 192  # J33: 462.22223 -> 192
$ws.Cells.Item(33, 11).Value = 1307.25  # K33: 1012.90908 -> 1307.25
$ws.Cells.Item(33, 12).Value = 1152  # L33: 2773.33338 -> 1152
$ws.Cells.Item(33, 13).Value = -1024.25  # M33: -729.9090800000001 -> -1024.25
$ws.Cells.Item(33, 14).Value = -1718  # N33: -3339.33338 -> -1718
$ws.Cells.Item(68, 8).Value = 511.64285  # H68: 472.0625 -> 511.64285
$ws.Cells.Item(68, 9).Value = 212.25  # I68: 206.5 -> 212.25
$ws.Cells.Item(68, 11).Value = 636.75  # K68: 619.5 -> 636.75
$ws.Cells.Item(68, 13).Value = 174.25  # M68: 191.5 -> 174.25
$ws.Cells.Item(71, 8).Value = 511.64285  # H71: 472.0625 -> 511.64285
$ws.Cells.Item(71, 9).Value = 212.25  # I71: 206.5 -> 212.25
$ws.Cells.Item(71, 11).Value = 1910.25  # K71: 1858.5 -> 1910.25
$ws.Cells.Item(71, 13).Value = 2145.75  # M71: 2197.5 -> 2145.75
$ws.Cells.Item(129, 8).Value = 3214.9375  # H129: 3439.5 -> 3214.9375
$ws.Cells.Item(129, 9).Value = 382.875  # I129: 395.2857 -> 382.875
$ws.Cells.Item(129, 10).Value = 6047  # J129: 6483.7144 -> 6047
$ws.Cells.Item(129, 11).Value = 1148.625  # K129: 1185.8571 -> 1148.625
$ws.Cells.Item(129, 12).Value = 18141  # L129: 19451.1432 -> 18141
$ws.Cells.Item(129, 13).Value = 3851.375  # M129: 3814.1429 -> 3851.375
$ws.Cells.Item(129, 14).Value = -28141  # N129: -29451.1432 -> -28141
$ws.Cells.Item(131, 8).Value = 1750.3334  # H131: 1811.625 -> 1750.3334
$ws.Cells.Item(131, 9).Value = 1020.375  # I131: 1076.1428 -> 1020.375
$ws.Cells.Item(131, 10).Value = 2334.3  # J131: 2383.6667 -> 2334.3
$ws.Cells.Item(131, 11).Value = 3061.125  # K131: 3228.4284 -> 3061.125
$ws.Cells.Item(131, 12).Value = 7002.900000000001  # L131: 7151.000100000001 -> 7002.900000000001
$ws.Cells.Item(131, 13).Value = 1978.875  # M131: 1811.5716 -> 1978.875
$ws.Cells.Item(131, 14).Value = -17082.9  # N131: -17231.0001 -> -17082.9
$ws.Cells.Item(135, 8).Value = 781.8889  # H135: 793.3 -> 781.8889
$ws.Cells.Item(135, 10).Value = 999.6667  # J135: 973.75 -> 999.6667
$ws.Cells.Item(135, 12).Value = 8997.0003  # L135: 8763.75 -> 8997.0003
$ws.Cells.Item(135, 14).Value = -14067.0003  # N135: -13833.75 -> -14067.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 678  # H102: 797.6667 -> 678
$ws.Cells.Item(102, 9).Value = 678  # I102: 797.6667 -> 678
$ws.Cells.Item(102, 11).Value = 678  # K102: 797.6667 -> 678
$ws.Cells.Item(102, 13).Value = 944  # M102: 824.3333 -> 944
$ws.Cells.Item(122, 8).Value = 1060.5385  # H122: 1138.8 -> 1060.5385
$ws.Cells.Item(122, 9).Value = 878.7  # I122: 923.5 -> 878.7
$ws.Cells.Item(122, 10).Value = 1666.6666  # J122: 2000 -> 1666.6666
$ws.Cells.Item(122, 11).Value = 2636.1  # K122: 2770.5 -> 2636.1
$ws.Cells.Item(122, 12).Value = 4999.9998  # L122: 6000 -> 4999.9998
$ws.Cells.Item(122, 13).Value = -186.1000000000004  # M122: -320.5 -> -186.1000000000004
$ws.Cells.Item(122, 14).Value = -9899.9998  # N122: -10900 -> -9899.9998

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 3505  # H22: 3705.875 -> 3505
$ws.Cells.Item(22, 9).Value = 3510.111  # I22: 3912.5 -> 3510.111
$ws.Cells.Item(22, 11).Value = 3510.111  # K22: 3912.5 -> 3510.111
$ws.Cells.Item(22, 13).Value = -3215.111  # M22: -3617.5 -> -3215.111
$ws.Cells.Item(27, 8).Value = 3505  # H27: 3705.875 -> 3505
$ws.Cells.Item(27, 9).Value = 3510.111  # I27: 3912.5 -> 3510.111
$ws.Cells.Item(27, 11).Value = 3510.111  # K27: 3912.5 -> 3510.111
$ws.Cells.Item(27, 13).Value = -3403.111  # M27: -3805.5 -> -3403.111
$ws.Cells.Item(40, 8).Value = 5657.3335  # H40: 5979.2104 -> 5657.3335
$ws.Cells.Item(40, 9).Value = 4156.6665  # I40: 4468.1 -> 4156.6665
$ws.Cells.Item(40, 11).Value = 4156.6665  # K40: 4468.1 -> 4156.6665
$ws.Cells.Item(40, 13).Value = -4020.6665  # M40: -4332.1 -> -4020.6665
$ws.Cells.Item(41, 8).Value = 30000  # H41: 0 -> 30000
$ws.Cells.Item(41, 10).Value = 30000  # J41: 0 -> 30000
$ws.Cells.Item(41, 12).Value = 30000  # L41: 0 -> 30000
$ws.Cells.Item(41, 14).Value = -30876  # N41: None -> -30876
$ws.Cells.Item(132, 8).Value = 3516.4644  # H132: 3272.9355 -> 3516.4644
$ws.Cells.Item(132, 9).Value = 2725.5386  # I132: 2402 -> 2725.5386
$ws.Cells.Item(132, 11).Value = 8176.6158  # K132: 7206 -> 8176.6158
$ws.Cells.Item(132, 13).Value = -5646.6158  # M132: -4676 -> -5646.6158

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(45, 8).Value = 14813  # H45: 15000 -> 14813
$ws.Cells.Item(45, 10).Value = 14626  # J45: 0 -> 14626
$ws.Cells.Item(45, 12).Value = 14626  # L45: 0 -> 14626
$ws.Cells.Item(45, 14).Value = -15608  # N45: None -> -15608
$ws.Cells.Item(126, 8).Value = 4804.1665  # H126: 5006.875 -> 4804.1665
$ws.Cells.Item(126, 9).Value = 4649.375  # I126: 4900 -> 4649.375
$ws.Cells.Item(126, 11).Value = 13948.125  # K126: 14700 -> 13948.125
$ws.Cells.Item(126, 13).Value = -11478.125  # M126: -12230 -> -11478.125

